$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast values and timestamps for rows 2-93 (new data window)
$ws.Cells.Item(2, 1).Value = 5280
$ws.Cells.Item(2, 2).Value = 45932
$ws.Cells.Item(3, 1).Value = 5250
$ws.Cells.Item(3, 2).Value = 45932.01041666666
$ws.Cells.Item(4, 1).Value = 5230
$ws.Cells.Item(4, 2).Value = 45932.02083333334
$ws.Cells.Item(5, 1).Value = 5200
$ws.Cells.Item(5, 2).Value = 45932.03125
$ws.Cells.Item(6, 1).Value = 5190
$ws.Cells.Item(6, 2).Value = 45932.04166666666
$ws.Cells.Item(7, 1).Value = 5170
$ws.Cells.Item(7, 2).Value = 45932.05208333334
$ws.Cells.Item(8, 1).Value = 5150
$ws.Cells.Item(8, 2).Value = 45932.0625
$ws.Cells.Item(9, 1).Value = 5150
$ws.Cells.Item(9, 2).Value = 45932.07291666666
$ws.Cells.Item(10, 1).Value = 5150
$ws.Cells.Item(10, 2).Value = 45932.08333333334
$ws.Cells.Item(11, 1).Value = 5150
$ws.Cells.Item(11, 2).Value = 45932.09375
$ws.Cells.Item(12, 1).Value = 5160
$ws.Cells.Item(12, 2).Value = 45932.10416666666
$ws.Cells.Item(13, 1).Value = 5170
$ws.Cells.Item(13, 2).Value = 45932.11458333334
$ws.Cells.Item(14, 1).Value = 5170
$ws.Cells.Item(14, 2).Value = 45932.125
$ws.Cells.Item(15, 1).Value = 5180
$ws.Cells.Item(15, 2).Value = 45932.13541666666
$ws.Cells.Item(16, 1).Value = 5190
$ws.Cells.Item(16, 2).Value = 45932.14583333334
$ws.Cells.Item(17, 1).Value = 5230
$ws.Cells.Item(17, 2).Value = 45932.15625
$ws.Cells.Item(18, 1).Value = 5290
$ws.Cells.Item(18, 2).Value = 45932.16666666666
$ws.Cells.Item(19, 1).Value = 5370
$ws.Cells.Item(19, 2).Value = 45932.17708333334
$ws.Cells.Item(20, 1).Value = 5470
$ws.Cells.Item(20, 2).Value = 45932.1875
$ws.Cells.Item(21, 1).Value = 5600
$ws.Cells.Item(21, 2).Value = 45932.19791666666
$ws.Cells.Item(22, 1).Value = 5740
$ws.Cells.Item(22, 2).Value = 45932.20833333334
$ws.Cells.Item(23, 1).Value = 5900
$ws.Cells.Item(23, 2).Value = 45932.21875
$ws.Cells.Item(24, 1).Value = 6070
$ws.Cells.Item(24, 2).Value = 45932.22916666666
$ws.Cells.Item(25, 1).Value = 6250
$ws.Cells.Item(25, 2).Value = 45932.23958333334
$ws.Cells.Item(26, 1).Value = 6420
$ws.Cells.Item(26, 2).Value = 45932.25
$ws.Cells.Item(27, 1).Value = 6590
$ws.Cells.Item(27, 2).Value = 45932.26041666666
$ws.Cells.Item(28, 1).Value = 6740
$ws.Cells.Item(28, 2).Value = 45932.27083333334
$ws.Cells.Item(29, 1).Value = 6880
$ws.Cells.Item(29, 2).Value = 45932.28125
$ws.Cells.Item(30, 1).Value = 6990
$ws.Cells.Item(30, 2).Value = 45932.29166666666
$ws.Cells.Item(31, 1).Value = 7090
$ws.Cells.Item(31, 2).Value = 45932.30208333334
$ws.Cells.Item(32, 1).Value = 7150
$ws.Cells.Item(32, 2).Value = 45932.3125
$ws.Cells.Item(33, 1).Value = 7150
$ws.Cells.Item(33, 2).Value = 45932.32291666666
$ws.Cells.Item(34, 1).Value = 7150
$ws.Cells.Item(34, 2).Value = 45932.33333333334
$ws.Cells.Item(35, 1).Value = 7150
$ws.Cells.Item(35, 2).Value = 45932.34375
$ws.Cells.Item(36, 1).Value = 7120
$ws.Cells.Item(36, 2).Value = 45932.35416666666
$ws.Cells.Item(37, 1).Value = 7050
$ws.Cells.Item(37, 2).Value = 45932.36458333334
$ws.Cells.Item(38, 1).Value = 6970
$ws.Cells.Item(38, 2).Value = 45932.375
$ws.Cells.Item(39, 1).Value = 6880
$ws.Cells.Item(39, 2).Value = 45932.38541666666
$ws.Cells.Item(40, 1).Value = 6790
$ws.Cells.Item(40, 2).Value = 45932.39583333334
$ws.Cells.Item(41, 1).Value = 6700
$ws.Cells.Item(41, 2).Value = 45932.40625
$ws.Cells.Item(42, 1).Value = 6620
$ws.Cells.Item(42, 2).Value = 45932.41666666666
$ws.Cells.Item(43, 1).Value = 6560
$ws.Cells.Item(43, 2).Value = 45932.42708333334
$ws.Cells.Item(44, 1).Value = 6500
$ws.Cells.Item(44, 2).Value = 45932.4375
$ws.Cells.Item(45, 1).Value = 6460
$ws.Cells.Item(45, 2).Value = 45932.44791666666
$ws.Cells.Item(46, 1).Value = 6440
$ws.Cells.Item(46, 2).Value = 45932.45833333334
$ws.Cells.Item(47, 1).Value = 6420
$ws.Cells.Item(47, 2).Value = 45932.46875
$ws.Cells.Item(48, 1).Value = 6410
$ws.Cells.Item(48, 2).Value = 45932.47916666666
$ws.Cells.Item(49, 1).Value = 6400
$ws.Cells.Item(49, 2).Value = 45932.48958333334
$ws.Cells.Item(50, 1).Value = 6400
$ws.Cells.Item(50, 2).Value = 45932.5
$ws.Cells.Item(51, 1).Value = 6390
$ws.Cells.Item(51, 2).Value = 45932.51041666666
$ws.Cells.Item(52, 1).Value = 6390
$ws.Cells.Item(52, 2).Value = 45932.52083333334
$ws.Cells.Item(53, 1).Value = 6380
$ws.Cells.Item(53, 2).Value = 45932.53125
$ws.Cells.Item(54, 1).Value = 6380
$ws.Cells.Item(54, 2).Value = 45932.54166666666
$ws.Cells.Item(55, 1).Value = 6370
$ws.Cells.Item(55, 2).Value = 45932.55208333334
$ws.Cells.Item(56, 1).Value = 6370
$ws.Cells.Item(56, 2).Value = 45932.5625
$ws.Cells.Item(57, 1).Value = 6370
$ws.Cells.Item(57, 2).Value = 45932.57291666666
$ws.Cells.Item(58, 1).Value = 6370
$ws.Cells.Item(58, 2).Value = 45932.58333333334
$ws.Cells.Item(59, 1).Value = 6380
$ws.Cells.Item(59, 2).Value = 45932.59375
$ws.Cells.Item(60, 1).Value = 6390
$ws.Cells.Item(60, 2).Value = 45932.60416666666
$ws.Cells.Item(61, 1).Value = 6410
$ws.Cells.Item(61, 2).Value = 45932.61458333334
$ws.Cells.Item(62, 1).Value = 6440
$ws.Cells.Item(62, 2).Value = 45932.625
$ws.Cells.Item(63, 1).Value = 6480
$ws.Cells.Item(63, 2).Value = 45932.63541666666
$ws.Cells.Item(64, 1).Value = 6520
$ws.Cells.Item(64, 2).Value = 45932.64583333334
$ws.Cells.Item(65, 1).Value = 6580
$ws.Cells.Item(65, 2).Value = 45932.65625
$ws.Cells.Item(66, 1).Value = 6630
$ws.Cells.Item(66, 2).Value = 45932.66666666666
$ws.Cells.Item(67, 1).Value = 6690
$ws.Cells.Item(67, 2).Value = 45932.67708333334
$ws.Cells.Item(68, 1).Value = 6760
$ws.Cells.Item(68, 2).Value = 45932.6875
$ws.Cells.Item(69, 1).Value = 6840
$ws.Cells.Item(69, 2).Value = 45932.69791666666
$ws.Cells.Item(70, 1).Value = 6930
$ws.Cells.Item(70, 2).Value = 45932.70833333334
$ws.Cells.Item(71, 1).Value = 7030
$ws.Cells.Item(71, 2).Value = 45932.71875
$ws.Cells.Item(72, 1).Value = 7150
$ws.Cells.Item(72, 2).Value = 45932.72916666666
$ws.Cells.Item(73, 1).Value = 7280
$ws.Cells.Item(73, 2).Value = 45932.73958333334
$ws.Cells.Item(74, 1).Value = 7410
$ws.Cells.Item(74, 2).Value = 45932.75
$ws.Cells.Item(75, 1).Value = 7520
$ws.Cells.Item(75, 2).Value = 45932.76041666666
$ws.Cells.Item(76, 1).Value = 7580
$ws.Cells.Item(76, 2).Value = 45932.77083333334
$ws.Cells.Item(77, 1).Value = 7620
$ws.Cells.Item(77, 2).Value = 45932.78125
$ws.Cells.Item(78, 1).Value = 7650
$ws.Cells.Item(78, 2).Value = 45932.79166666666
$ws.Cells.Item(79, 1).Value = 7600
$ws.Cells.Item(79, 2).Value = 45932.80208333334
$ws.Cells.Item(80, 1).Value = 7500
$ws.Cells.Item(80, 2).Value = 45932.8125
$ws.Cells.Item(81, 1).Value = 7410
$ws.Cells.Item(81, 2).Value = 45932.82291666666
$ws.Cells.Item(82, 1).Value = 7300
$ws.Cells.Item(82, 2).Value = 45932.83333333334
$ws.Cells.Item(83, 1).Value = 7160
$ws.Cells.Item(83, 2).Value = 45932.84375
$ws.Cells.Item(84, 1).Value = 7060
$ws.Cells.Item(84, 2).Value = 45932.85416666666
$ws.Cells.Item(85, 1).Value = 6930
$ws.Cells.Item(85, 2).Value = 45932.86458333334
$ws.Cells.Item(86, 1).Value = 6760
$ws.Cells.Item(86, 2).Value = 45932.875
$ws.Cells.Item(87, 1).Value = 6640
$ws.Cells.Item(87, 2).Value = 45932.88541666666
$ws.Cells.Item(88, 1).Value = 6480
$ws.Cells.Item(88, 2).Value = 45932.89583333334
$ws.Cells.Item(89, 1).Value = 6310
$ws.Cells.Item(89, 2).Value = 45932.90625
$ws.Cells.Item(90, 1).Value = 6100
$ws.Cells.Item(90, 2).Value = 45932.91666666666
$ws.Cells.Item(91, 1).Value = 5910
$ws.Cells.Item(91, 2).Value = 45932.92708333334
$ws.Cells.Item(92, 1).Value = 5830
$ws.Cells.Item(92, 2).Value = 45932.9375
$ws.Cells.Item(93, 1).Value = 5720
$ws.Cells.Item(93, 2).Value = 45932.94791666666

# Remove now-unused trailing rows 94-97 (data window shortened)
$ws.Range("A94:B97").EntireRow.Delete() | Out-Null

Write-Host ("Final UsedRange rows: " + $ws.UsedRange.Rows.Count())
Write-Host ("Dimension check A93: " + $ws.Cells.Item(93,1).Value())
